$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.795.76'
$ws.Range("E2").Value = '  -2.39%  '

$ws.Range("D3").Value = '3.271.66'
$ws.Range("E3").Value = '  -1.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.08'
$ws.Range("E5").Value = '  -1.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.27'
$ws.Range("E6").Value = '  -4.98%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  +4.03%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("E9").Value = '  -3.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.69'
$ws.Range("E10").Value = '  +0.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.398'
$ws.Range("E11").Value = '  -2.71%  '

$ws.Range("D12").Value = '3.840.37'
$ws.Range("E12").Value = '  -1.21%  '

$ws.Range("E13").Value = '  -3.79%  '

$ws.Range("D14").Value = '65.883.36'
$ws.Range("E14").Value = '  -2.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.40'
$ws.Range("E15").Value = '  -3.89%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000162'
$ws.Range("E16").Value = '  -3.30%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.245.51'
$ws.Range("E17").Value = '  -1.65%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '436.18'
$ws.Range("E18").Value = '  -1.96%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.54'
$ws.Range("E19").Value = '  -2.81%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.13'
$ws.Range("E20").Value = '  -3.39%  '

$ws.Range("E21").Value = '  -4.98%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.18'
$ws.Range("E22").Value = '  -2.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("D24").Value = '3.420.93'
$ws.Range("E24").Value = '  -1.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.505'
$ws.Range("E25").Value = '  -2.62%  '

$ws.Range("E26").Value = '  +3.42%  '

$ws.Range("E27").Value = '  -5.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.84'
$ws.Range("E28").Value = '  -2.44%  '

$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("E30").Value = '  -2.67%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.23'
$ws.Range("E31").Value = '  -3.16%  '

$ws.Range("E32").Value = '  +0.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.13'
$ws.Range("E33").Value = '  -4.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.58'
$ws.Range("E34").Value = '  -3.49%  '

$ws.Range("E35").Value = '  -5.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.03'
$ws.Range("E36").Value = '  -1.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.43'
$ws.Range("E37").Value = '  -5.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.66'
$ws.Range("E38").Value = '  -2.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.78'
$ws.Range("E39").Value = '  -4.32%  '

$ws.Range("D40").Value = '2.757.96'
$ws.Range("E40").Value = '  -0.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.774'
$ws.Range("E41").Value = '  -2.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.31'
$ws.Range("E42").Value = '  -3.84%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.24'
$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.02'
$ws.Range("E44").Value = '  -3.98%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0655'
$ws.Range("E45").Value = '  -2.93%  '

$ws.Range("E46").Value = '  -6.15%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '317.52'
$ws.Range("E47").Value = '  -2.86%  '

$ws.Range("E48").Value = '  -6.77%  '

$ws.Range("E49").Value = '  -3.14%  '

$ws.Range("E50").Value = '  +2.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  +0.01%  '

